# Commit: "updated with Cornell and Penn data."
# Fix the "#NN- Name" player labels so there is a space before the dash,
# e.g. "#35- W. Cheek" -> "#35 - W. Cheek" (column H, "player").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H11").Value = "#35 - W. Cheek"
$ws.Range("H14").Value = "#48 - K. Salvatore"
$ws.Range("H17").Value = "#9 - S. Evans"
$ws.Range("H19").Value = "#35 - W. Cheek"
$ws.Range("H23").Value = "#24 - J. Lang"
$ws.Range("H24").Value = "#24 - J. Lang"
$ws.Range("H27").Value = "#26 - N. Andrews"

# Update the sheet's current selection to match the saved view state.
$ws.Range("H28").Select()
